# ScouseTom_New Pick and Place - add R8..R11, renumber query/name, refresh selection.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The refreshed .mnt import picked up four extra resistors (R8-R11) ahead of
# the existing U1-U5 rows. Insert four blank rows at the old row 48 so the
# U1..U5 block slides down to 52..56, then populate the new rows.
$ws.Rows("48:51").Insert()

$ws.Range("A48").Value = "R8"
$ws.Range("B48").Value = -2.65
$ws.Range("C48").Value = 39.15
$ws.Range("D48").Value = 270
$ws.Range("E48").Value = "10k"
$ws.Range("F48").Value = "R0805"
$ws.Range("G48").Value = "MC01W0805110K"
$ws.Range("H48").Value = "Top"

$ws.Range("A49").Value = "R9"
$ws.Range("B49").Value = -2.25
$ws.Range("C49").Value = 33.05
$ws.Range("D49").Value = 270
$ws.Range("E49").Value = "10k"
$ws.Range("F49").Value = "R0805"
$ws.Range("G49").Value = "MC01W0805110K"
$ws.Range("H49").Value = "Top"

$ws.Range("A50").Value = "R10"
$ws.Range("B50").Value = 100.35
$ws.Range("C50").Value = 48.15
$ws.Range("D50").Value = 90
$ws.Range("E50").Value = "10k"
$ws.Range("F50").Value = "R0805"
$ws.Range("G50").Value = "MC01W0805110K"
$ws.Range("H50").Value = "Top"

$ws.Range("A51").Value = "R11"
$ws.Range("B51").Value = 100.2
$ws.Range("C51").Value = 35.1
$ws.Range("D51").Value = 90
$ws.Range("E51").Value = "10k"
$ws.Range("F51").Value = "R0805"
$ws.Range("G51").Value = "MC01W0805110K"
$ws.Range("H51").Value = "Top"

# The query/defined-name pair gets bumped to _1 now that the range covers
# the four extra rows (A2:F52 -> A2:F56).
$n = $wb.Names.Item(1)
$n.Name = "ScouseTom_New_1"
$n.RefersTo = "=Sheet1!`$A`$2:`$F`$56"

# Leave the selection where the user last clicked after reviewing the new rows.
[void]$ws.Range("K11").Select()
